$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set text values in the exact order that yields the target shared-string table order ---
$ws.Range("C7").Value = ' Want to hear what we\''re up to?'
$ws.Range("A7").Value = 'SCRIPT/T01P02A/um1202.ssb'
$ws.Range("C8").Value = ' We\''re on the hunt for the [CS:I]Golden\nMask[CR].'
$ws.Range("C9").Value = ' It\''s an ancient mask from a\nlong-lost civilization.'
$ws.Range("C10").Value = ' Would you two know anything\nabout it?'
$ws.Range("D7").Value = ' Хочешь знать, что мы ищем?'
$ws.Range("D8").Value = ' Мы охотимся за [CS:I]Золотой Маской[CR].'
$ws.Range("D9").Value = ' Это древняя маска давно\nисчезнувшей цивилизации.'
$ws.Range("D10").Value = ' Вы ребята, знаете что-нибудь\nо ней?'
$ws.Range("E7").Value = ' Öïœåšû èîàóû, œóï íú éþåí?'
$ws.Range("E8").Value = ' Íú ïöïóéíòÿ èà [CS:I]Èïìïóïê Íàòëïê[CR].'
$ws.Range("E9").Value = ' Üóï äñåâîÿÿ íàòëà äàâîï\néòœåèîôâšåê øéâéìéèàøéé.'
$ws.Range("E10").Value = ' Âú ñåáÿóà, èîàåóå œóï-îéáôäû\nï îåê?'
$ws.Range("A8").Value = 'SCRIPT/T01P02A/um1204.ssb'
$ws.Range("A11").Value = 'SCRIPT/P01P04A/um1405.ssb'
$ws.Range("C11").Value = ' I won on the first [CS:I]Prize Ticket[CR]\nI got.'
$ws.Range("C12").Value = ' It must be beginner\''s luck.\nHa ha ha!'
$ws.Range("D11").Value = ' Я получил выигрыш с самого\nпервого [CS:I]Призового Билета[CR].'
$ws.Range("D12").Value = ' Должно быть, это удача новичка.\nХа-ха-ха!'
$ws.Range("E11").Value = ' Ÿ ðïìôœéì âúéãñúš ò òàíïãï\nðåñâïãï [CS:I]Ðñéèïâïãï Áéìåóà[CR].'
$ws.Range("E12").Value = ' Äïìçîï áúóû, üóï ôäàœà îïâéœëà.\nÖà-öà-öà!'
$ws.Range("A12").Value = 'SCRIPT/P01P04A/um1505.ssb'
$ws.Range("A13").Value = 'SCRIPT/G01P03A/um1614.ssb'
$ws.Range("C13").Value = ' No word from the great\n[CS:N]Dusknoir[CR] yet?'
$ws.Range("D13").Value = ' От великого [CS:N]Даскнуара[CR] ещё не было\nвестей?'
$ws.Range("E13").Value = ' Ïó âåìéëïãï [CS:N]Äàòëîôàñà[CR] åþæ îå áúìï\nâåòóåê?'

# --- Set the numeric "line number" column B ---
$ws.Range("B7").Value = 155
$ws.Range("B8").Value = 158
$ws.Range("B9").Value = 161
$ws.Range("B10").Value = 164
$ws.Range("B11").Value = 123
$ws.Range("B12").Value = 126
$ws.Range("B13").Value = 104

# --- Re-apply the alternating (bottom-border) style to rows 6, 10 and 12, like rows 3 in the original table ---
$ws.Range("A3:E3").Copy()
$ws.Range("A6:E6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A10:E10").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A12:E12").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# --- Row heights ---
$ws.Rows.Item(7).RowHeight = 43.2
$ws.Rows.Item(8).RowHeight = 43.2
$ws.Rows.Item(9).RowHeight = 21.6
$ws.Rows.Item(10).RowHeight = 21.6
$ws.Rows.Item(11).RowHeight = 43.2
$ws.Rows.Item(12).RowHeight = 41.4
$ws.Rows.Item(13).RowHeight = 43.2

# --- Selection / scroll position, matching the view state recorded in the workbook ---
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D13").Select()
